# Localization status report refresh ("Generate Report for Archive"):
# the two handed-off docs have moved from "Ready for handoff" into
# "In Translation" -- update the status everywhere it is shown
# (Overview's per-language columns + each language sheet's Status column)
# and shrink the now-narrower Status columns to fit the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) / de-de (col F) status for both rows ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Narrow the two status columns now that the text is shorter.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status column (C) for both data rows ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status column (C) for both data rows ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = 12.5
